$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Friday (2018/3/9) weekly-plan block, rows 26-28 (C/D content),
# mirroring the structure used for the other days above it.
$ws.Range("A26").Value = "2018/3/9`n周五"
$ws.Range("B26").Value = "上午"

$ws.Range("C26").Value = "修改完成sf_cap.c和sf_cap.h中的错误"
$ws.Range("D26").Value = "9:20 - 10:20"

$ws.Range("C27").Value = "修改完成qos中添加action的接口"
$ws.Range("D27").Value = "10:30 - 11:40"

$ws.Range("C28").Value = "背单词40个"

# Update the active selection to match the final workbook state.
$ws.Range("H26:H31").Select()
